# Refresh the cryptos list (price + 1h volume change) per the
# "Updated cryptos list ... with GitHub Actions" commit.
#
# D (Price) / E (Volume(1h)) cells are plain text in the source workbook
# (e.g. "58.602.68", "1.00", "  +3.95%  "), so for price cells that would
# otherwise be auto-parsed as a clean number (losing trailing zeros, e.g.
# "1.00" -> 1), we force the cell to Text format first so the literal
# string is preserved, matching the original inlineStr content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.602.68'
$ws.Range("E2").Value = '  +3.95%  '

$ws.Range("D3").Value = '3.302.59'
$ws.Range("E3").Value = '  +2.21%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '400.35'
$ws.Range("E5").Value = '  +0.62%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '109.68'
$ws.Range("E6").Value = '  -1.63%  '

$ws.Range("E7").Value = '  +5.49%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.635'
$ws.Range("E9").Value = '  +2.28%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.86'
$ws.Range("E10").Value = '  +1.32%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0980'
$ws.Range("E11").Value = '  +7.15%  '

$ws.Range("E12").Value = '  +1.38%  '

$ws.Range("D13").Value = '3.816.13'
$ws.Range("E13").Value = '  +2.04%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.37'
$ws.Range("E14").Value = '  +2.79%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '19.28'
$ws.Range("E15").Value = '  +0.92%  '

$ws.Range("D16").Value = '3.294.61'
$ws.Range("E16").Value = '  +2.19%  '

$ws.Range("E17").Value = '  -0.49%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.97'
$ws.Range("E18").Value = '  +0.41%  '

$ws.Range("D19").Value = '58.374.98'
$ws.Range("E19").Value = '  +3.81%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.32'
$ws.Range("E20").Value = '  -0.68%  '

$ws.Range("E21").Value = '  +6.19%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '12.96'
$ws.Range("E22").Value = '  -0.55%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '303.90'
$ws.Range("E23").Value = '  +1.60%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '74.81'

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.19'
$ws.Range("E25").Value = '  -0.93%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '28.36'
$ws.Range("E26").Value = '  +0.65%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '4.43'
$ws.Range("E27").Value = '  +1.24%  '

$ws.Range("E28").Value = '  -3.54%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.41'
$ws.Range("E29").Value = '  -0.68%  '

$ws.Range("E30").Value = '  -1.82%  '

$ws.Range("E31").Value = '  -0.40%  '

$ws.Range("E32").Value = '  +2.43%  '

$ws.Range("E33").Value = '  +1.85%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '40.80'
$ws.Range("E34").Value = '  +11.19%  '

$ws.Range("E35").Value = '  +6.41%  '

$ws.Range("E36").Value = '  -2.16%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '51.90'
$ws.Range("E37").Value = '  +0.98%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.30'
$ws.Range("E38").Value = '  +5.57%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.00'
$ws.Range("E39").Value = '  +0.12%  '

$ws.Range("E40").Value = '  -1.12%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '137.84'
$ws.Range("E41").Value = '  +0.32%  '

$ws.Range("E42").Value = '  +2.74%  '

$ws.Range("B43").Value = 'NEARProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.95'
$ws.Range("E43").Value = '  -1.43%  '

$ws.Range("B44").Value = 'ARBITRUM'
$ws.Range("C44").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.88'
$ws.Range("E44").Value = '  -2.07%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '16.87'
$ws.Range("E45").Value = '  -3.23%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.280'
$ws.Range("E46").Value = '  -1.98%  '

$ws.Range("E47").Value = '  +9.31%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '22.59'
$ws.Range("E48").Value = '  +1.65%  '

$ws.Range("D49").Value = '2.165.87'
$ws.Range("E49").Value = '  +1.73%  '

$ws.Range("E50").Value = '  -0.49%  '

$ws.Range("E51").Value = '  -13.90%  '
